# Replace the two M2Doc field codes ( {m:userdoc 'zone1'} and {m:enduserdoc} )
# that are currently stored as real Word fields (fldChar begin/instrText/fldChar end)
# with plain literal text runs containing the braces, matching the
# TokenIteratorFieldRewriterSplit output format.

$d = $word.ActiveDocument

$pkgHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>'
$pkgFooter = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- First field: { m:userdoc 'zone1' } -> split across 4 runs ---
$field1 = $d.Fields.Item(1)
$insertAt1 = $field1.Code.Start - 1
$field1.Delete()

$body1 = '<w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:userdoc ''zone1''</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>'
$xml1 = $pkgHeader + $body1 + $pkgFooter
$range1 = $d.Range($insertAt1, $insertAt1)
$range1.InsertXML($xml1)

# --- Second field: { m:enduserdoc } -> split into 2 runs around the _GoBack bookmark ---
$field2 = $d.Fields.Item(1)
$insertAt2 = $field2.Code.Start - 1
$field2.Delete()

$body2 = '<w:p><w:r><w:t>{m:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">enduserdoc}</w:t></w:r></w:p>'
$xml2 = $pkgHeader + $body2 + $pkgFooter
$range2 = $d.Range($insertAt2, $insertAt2)
$range2.InsertXML($xml2)
